$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last fully-populated data row (48) down into the
# three new rows (51-53) so the new cells pick up the same styles (s="1" for
# A/B, s="2" for C:U) used throughout the table.
$ws.Range("A48:U48").Copy()
$ws.Range("A51:U53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 51: Gatchina MR, 2019 ---
$ws.Range("A51").Value = "Гатчинский МР"
$ws.Range("B51").Value = 2019
$ws.Range("C51").Value = 238.018
$ws.Range("D51").Value = 42.329000000000001
$ws.Range("E51").Value = 1035
$ws.Range("F51").Value = 47095.7
$ws.Range("G51").Value = "???"
$ws.Range("H51").Value = "???"
$ws.Range("I51").Value = "???"
$ws.Range("J51").Value = "???"
$ws.Range("K51").Value = "???"
$ws.Range("L51").Value = 8534.4
$ws.Range("M51").Value = "???"
$ws.Range("N51").Value = 8877
$ws.Range("O51").Value = 98966.3
$ws.Range("P51").Value = 10252.700000000001
$ws.Range("Q51").Value = 88.7
$ws.Range("R51").Value = 29
$ws.Range("S51").Formula = "=22883164.3/1000"
$ws.Range("T51").Formula = "=353488/1000"
$ws.Range("U51").Value = -3761

# --- Row 52: Gatchina MR, 2018 ---
$ws.Range("A52").Value = "Гатчинский МР"
$ws.Range("B52").Value = 2018
$ws.Range("C52").Value = 243.17
$ws.Range("D52").Value = 43.061
$ws.Range("E52").Value = 836
$ws.Range("F52").Value = 43057.3
$ws.Range("G52").Value = "???"
$ws.Range("H52").Value = "???"
$ws.Range("I52").Value = "???"
$ws.Range("J52").Value = "???"
$ws.Range("K52").Value = "???"
$ws.Range("L52").Value = 19293.900000000001
$ws.Range("M52").Value = "???"
$ws.Range("N52").Value = 8284
$ws.Range("O52").Value = 81900
$ws.Range("P52").Value = 8500
$ws.Range("Q52").Value = 142.5
$ws.Range("R52").Value = 646
$ws.Range("S52").Formula = "=16022223.4/1000"
$ws.Range("T52").Formula = "=171538.4/1000"
$ws.Range("U52").Value = 339

# --- Row 53: Gatchina MR, 2017 ---
$ws.Range("A53").Value = "Гатчинский МР"
$ws.Range("B53").Value = 2017
$ws.Range("C53").Value = 244.25800000000001
$ws.Range("D53").Value = 37.133000000000003
$ws.Range("E53").Value = 1006
$ws.Range("F53").Value = 41288
$ws.Range("G53").Value = "???"
$ws.Range("H53").Value = "???"
$ws.Range("I53").Value = "???"
$ws.Range("J53").Value = "???"
$ws.Range("K53").Value = "???"
$ws.Range("L53").Value = 6649.3
$ws.Range("M53").Value = "???"
$ws.Range("N53").Value = 7663
$ws.Range("O53").Value = 63344.5
$ws.Range("P53").Value = 2696.5
$ws.Range("Q53").Value = 165.1
$ws.Range("R53").Value = 1535
$ws.Range("S53").Formula = "=14877213.8/1000"
$ws.Range("T53").Formula = "=132601/1000"
$ws.Range("U53").Value = -63

# Reflect the updated scroll position / selection recorded in the sheet view.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("O55").Select() | Out-Null
